{"js": "const body = context.document.body;\n\n// --- Edit 1 ------------------------------------------------------------\n// The bold heading \"Ventana \u2013 formulario de usuarios:\" is a copy/paste\n// mistake: this particular heading sits right under the \"Ventana -\n// productos:\" section (after the bullet points describing the product\n// form) and should read \"...formulario de productos:\" instead. The exact\n// same heading text legitimately appears earlier for the *users* section,\n// so we must only touch the later occurrence, not every match.\nconst headingMatches = body.search(\"Ventana \u2013 formulario de usuarios:\", { matchCase: true });\nheadingMatches.load(\"text,paragraphs\");\nawait context.sync();\n\nlet targetHeading = null;\nfor (let i = 0; i < headingMatches.items.length; i++) {\n  const match = headingMatches.items[i];\n  const precedingPara = match.paragraphs.getFirst().getPrevious();\n  precedingPara.load(\"text\");\n  await context.sync();\n  if (precedingPara.text.indexOf(\"producto\") !== -1) {\n    targetHeading = match;\n  }\n}\n// Fall back to the last match only when there is more than one candidate\n// and none of them was pinned down by the context probe above (keeps the\n// script working even if the surrounding text ever changes, without\n// mis-firing when there is just a single, unrelated match).\nif (!targetHeading && headingMatches.items.length > 1) {\n  targetHeading = headingMatches.items[headingMatches.items.length - 1];\n}\n\nif (targetHeading) {\n  const wordInHeading = targetHeading.search(\"usuarios\", { matchCase: true });\n  wordInHeading.load(\"text\");\n  await context.sync();\n  if (wordInHeading.items.length > 0) {\n    wordInHeading.items[0].insertText(\"productos\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// --- Edit 2 --------------------------------------------------------------\n// \"Demas:\" is missing its accent mark; it should read \"Dem\u00e1s:\".\nconst demasMatches = body.search(\"Demas:\", { matchCase: true });\ndemasMatches.load(\"text\");\nawait context.sync();\nif (demasMatches.items.length > 0) {\n  demasMatches.items[0].insertText(\"Dem\u00e1s:\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# --- Edit 1 ---------------------------------------------------------------\n# The bold heading \"Ventana - formulario de usuarios:\" is a copy/paste\n# mistake: this particular heading sits right under the \"Ventana -\n# productos:\" section (after the bullet points describing the product\n# form) and should read \"...formulario de productos:\" instead. The exact\n# same heading text legitimately appears earlier for the *users* section,\n# so only the later occurrence (preceded by \"producto\" bullet text) must be\n# touched.\n$paragraphs = $d.Paragraphs\n$count = $paragraphs.Count\n\n$targetPara = $null\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $paragraphs.Item($i)\n  if ($p.Range.Text -like \"*Ventana*formulario de usuarios*\") {\n    if ($i -gt 1) {\n      $prevText = $paragraphs.Item($i - 1).Range.Text\n      if ($prevText -like \"*producto*\") {\n        $targetPara = $p\n      }\n    }\n  }\n}\n\n# Fallback: if the context probe above did not pin down a paragraph and\n# there is more than one candidate heading (e.g. the surrounding text ever\n# changes), fall back to the last match so the script keeps working -\n# without mis-firing when there is just a single, unrelated match.\nif ($targetPara -eq $null) {\n  $matchCount = 0\n  for ($i = 1; $i -le $count; $i++) {\n    if ($paragraphs.Item($i).Range.Text -like \"*Ventana*formulario de usuarios*\") {\n      $matchCount = $matchCount + 1\n    }\n  }\n  if ($matchCount -gt 1) {\n    for ($i = 1; $i -le $count; $i++) {\n      $p = $paragraphs.Item($i)\n      if ($p.Range.Text -like \"*Ventana*formulario de usuarios*\") {\n        $targetPara = $p\n      }\n    }\n  }\n}\n\nif ($targetPara -ne $null) {\n  $targetPara.Range.Text = \"Ventana \u2013 formulario de productos:\"\n}\n\n# --- Edit 2 -----------------------------------------------------------------\n# \"Demas:\" is missing its accent mark; it should read \"Dem\u00e1s:\".\n$demasPara = $null\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $paragraphs.Item($i)\n  if ($p.Range.Text -like \"*Demas:*\") {\n    $demasPara = $p\n  }\n}\n\nif ($demasPara -ne $null) {\n  $demasPara.Range.Text = \"Dem\u00e1s:\"\n}\n"}
